# edit.ps1 - apply the "tweaking the build system" revision to readme.docx
#
# Summary of changes:
#  1. "January 19, 2015" -> "January 20, 2015"
#  2. Move the hidden "_GoBack" bookmark from the empty paragraph after the
#     date line to inside the "Highlights of this release..." sentence
#     (right after "...this relea").
#  3. "Installers for Ubuntu and Debian based Linux distributions  (new)"
#     -> "Finally, added installers for Ubuntu and Debian based Linux
#         distributions  overhauling the build system(new)"
#  4. "Optimized stack operation for floats (new)" -> "Optimized stack
#     operations for floating point operations (new)"
#  5. "Added missing check for array assignment to array element
#      (bug/medium)" -> "Added missing compiler check for array assignment
#      to array element (bug/medium)"
#  6. Drop the stray <w:lastRenderedPageBreak/> before "Compile/execute for
#     code that has library dependencies:"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Date: "January 19" -> "January 20"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("January 19", $true, $false, $false, $false, $false, `
    $true, 1, $false, "January 20", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Relocate the "_GoBack" bookmark
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$r = $d.Content
$found = $r.Find.Execute("Highlights of this release include the following:", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $r.Start + 24   # right after "Highlights of this relea"
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------
# 3. "Installers for ..." -> "Finally, added installers for ..." plus a
#    new trailing clause about overhauling the build system.
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Installers for ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Finally, added installers for ", 2)

$r = $d.Content
$found = $r.Find.Execute("Finally, added installers for", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $afterInstallers = $r.End
    $rest = $d.Range($afterInstallers, $d.Content.End)
    $found2 = $rest.Find.Execute("based Linux distributions", $true, $false, `
        $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $secondSpace = $d.Range($rest.End + 1, $rest.End + 2)
        $secondSpace.Text = " overhauling the build system"
    }
}

# ---------------------------------------------------------------------
# 4. "Optimized stack operation for floats " -> "Optimized stack
#    operations for floating point operations "
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Optimized stack operation for floats ", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    "Optimized stack operations for floating point operations ", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. "Added missing check for array assignment to array element " ->
#    "Added missing compiler check for array assignment to array element "
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Added missing check for array assignment to array element ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Added missing compiler check for array assignment to array element ", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Drop the lastRenderedPageBreak marker before the "Compile/execute..."
#    run (re-typing the run text clears the stale render-break marker).
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Compile/execute for code that has library dependencies:", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.InsertBefore("")
}

Write-Output "edit.ps1 completed"
